$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transformer")

# The "transformer" sheet was missing the "b" (susceptance) column that the
# sibling "branch" sheet already has between "x" and "ShortTermRating".
# Insert a new column H, give it the header "b", and fill the data rows
# with 0 (matching the format used elsewhere in the workbook).
$ws.Columns("H").Insert()
$ws.Range("H1").Value = "b"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0

# Make "transformer" the active sheet/selection (was "generator" before).
$null = $ws.Range("H5").Select()
